$wb = $excel.ActiveWorkbook

# "zh-cn" sheet: update Correspond Handoff Datetime (E2) and
# Correspond Handback DateTime (H2)
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E2").Value = "2016-03-31 07:42:08"
$wsZhCn.Range("H2").Value = "2016-03-31 07:42:52"

# "de-de" sheet: update Correspond Handoff Datetime (E2) and
# Correspond Handback DateTime (H2)
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E2").Value = "2016-03-31 07:42:18"
$wsDeDe.Range("H2").Value = "2016-03-31 07:43:08"
